# "Actualización desde MV -datos-"
# Appends 5 new daily rows (07-09-2021 .. 13-09-2021) to the bottom of the
# liquidity injection/drain table on Sheet1 (rows 176-180, columns A:K).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: date label (col A, stored as text) followed by the 10 numeric
# values for columns B..K.
$newRows = @(
    @("07-09-2021", 38, 7, 9, -2, 0, -2, -8, -38, -11, 7),
    @("08-09-2021", 38, 7, 9, -2, 0, -2, -7, -39, -10, 7),
    @("09-09-2021", 37, 7, 9, -2, 0, -2, -3, -40, -12, 7),
    @("10-09-2021", 37, 7, 9, -2, 0, -2, -1, -39, -14, 7),
    @("13-09-2021", 37, 7, 9, -2, 0, -2, -2, -39, -13, 7)
)

$startRow = 176

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $values = $newRows[$i]

    # Column A holds a date-formatted label like the "Serie" cells above it
    # (e.g. "06-09-2021"), which must stay plain text. Assigning the literal
    # string directly gets auto-converted into a real date serial by Excel's
    # input parser, so instead build it via a text formula (forces a string
    # result) and then paste-special "values only" over itself to collapse
    # it down to a literal shared-string cell with no residual formula and
    # no number-format/style override.
    $cell = $ws.Cells.Item($row, 1)
    $cell.Formula = "=""" + $values[0] + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)

    for ($c = 1; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
